$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.076.94"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "1.790.03"
$ws.Range("E3").Value = "  -0.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "

# Row 6
$ws.Range("E6").Value = "  -1.18%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "

# Row 9
$ws.Range("E9").Value = "  +3.77%  "

# Row 10
$ws.Range("E10").Value = "  -3.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.12%  "

# Row 12
$ws.Range("D12").Value = "2.047.80"
$ws.Range("E12").Value = "  +0.00%  "

# Row 13
$ws.Range("E13").Value = "  +4.08%  "

# Row 14
$ws.Range("D14").Value = "1.787.90"
$ws.Range("E14").Value = "  -0.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.625"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "

# Row 16
$ws.Range("D16").Value = "34.063.09"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17
$ws.Range("E17").Value = "  +0.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "

# Row 20
$ws.Range("E20").Value = "  -1.27%  "

# Row 21
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("E23").Value = "  +0.38%  "

# Row 24
$ws.Range("E24").Value = "  -2.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "

# Row 26
$ws.Range("E26").Value = "  +1.68%  "

# Row 27
$ws.Range("E27").Value = "  -0.51%  "

# Row 28
$ws.Range("E28").Value = "  +0.90%  "

# Row 29
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
$ws.Range("E30").Value = "  +2.99%  "

# Row 31
$ws.Range("E31").Value = "  -0.52%  "

# Row 32
$ws.Range("E32").Value = "  -0.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.92%  "

# Row 34
$ws.Range("E34").Value = "  +1.25%  "

# Row 35
$ws.Range("D35").Value = "1.408.34"
$ws.Range("E35").Value = "  +1.58%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.655"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.31%  "

# Row 38
$ws.Range("E38").Value = "  +2.03%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "

# Row 41
$ws.Range("E41").Value = "  +0.48%  "

# Row 42
$ws.Range("E42").Value = "  +0.35%  "

# Row 43
$ws.Range("E43").Value = "  -2.67%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.51%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0140"
$ws.Range("E45").Value = "  +2.58%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.02%  "

# Row 47
$ws.Range("E47").Value = "  +1.65%  "

# Row 48
$ws.Range("E48").Value = "  +2.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "1.949.54"
$ws.Range("E50").Value = "  -0.34%  "

# Row 51
$ws.Range("E51").Value = "  -0.16%  "
